$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 15. This shifts the
# existing rows 15-17 down to 16-18 (their data is preserved as-is),
# matching the diff where row15->16, row16->17, row17->18 (new row 18).
$ws.Rows("15:15").Insert()

# Populate the newly-inserted row 15 with this week's record.
$ws.Range("A15").Value2 = 8
$ws.Range("B15").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C15").Value2 = "Coquimbo"
$ws.Range("D15").Value2 = 45166
$ws.Range("E15").Value2 = 4
$ws.Range("F15").Value2 = 100114002
$ws.Range("G15").Value2 = "Camote"
$ws.Range("H15").Value2 = "Sin especificar"
$ws.Range("I15").Value2 = "Primera"
$ws.Range("J15").Value2 = 200
$ws.Range("K15").Value2 = 16000
$ws.Range("L15").Value2 = 17000
$ws.Range("M15").Value2 = 16500
$ws.Range("N15").Value2 = "$/malla 18 kilos"
$ws.Range("O15").Value2 = "Perú"
$ws.Range("P15").Value2 = 917
$ws.Range("Q15").Value2 = 18
$ws.Range("R15").Value2 = "Hortaliza"
